$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @("D2", "311.77"),
    @("E2", "0.95%"),
    @("E3", "-2.84%"),
    @("D4", "5.103"),
    @("E4", "-2.55%"),
    @("D5", "0.07570"),
    @("E5", "-2.33%"),
    @("D6", "4.298"),
    @("E6", "-0.50%"),
    @("D7", "1.667"),
    @("E7", "2.92%"),
    @("D8", "0.9291"),
    @("E8", "1.04%"),
    @("E9", "-0.21%"),
    @("D10", "0.1210"),
    @("E10", "-2.49%"),
    @("D11", "0.1813"),
    @("E11", "-1.40%"),
    @("D12", "0.09035"),
    @("E12", "-0.75%"),
    @("D13", "0.04149"),
    @("E13", "-3.29%"),
    @("E14", "0.15%"),
    @("D15", "0.001285"),
    @("E15", "1.77%"),
    @("D16", "0.005943"),
    @("E16", "2.17%"),
    @("E18", "-0.09%"),
    @("E19", "0.59%"),
    @("D20", "7.652"),
    @("E20", "5.97%"),
    @("E21", "-2.24%"),
    @("D22", "0.2810"),
    @("E22", "-2.87%"),
    @("D23", "0.04017"),
    @("E23", "-1.29%"),
    @("D24", "0.001266"),
    @("E24", "0.30%"),
    @("D25", "0.004075"),
    @("E25", "-1.55%"),
    @("D38", "0.02422"),
    @("E38", "-1.11%"),
    @("D39", "0.05145"),
    @("E39", "-2.81%"),
    @("D40", "0.007716"),
    @("E40", "-1.74%"),
    @("D41", "0.1297"),
    @("E41", "-1.21%"),
    @("D42", "0.007645"),
    @("E42", "12.06%"),
    @("E43", "14.28%"),
    @("D44", "0.008066"),
    @("E44", "3.55%"),
    @("D45", "0.3112"),
    @("E45", "1.49%"),
    @("D46", "0.00006605"),
    @("E46", "-0.94%"),
    @("E47", "-0.14%"),
    @("D48", "1.572"),
    @("E48", "825.57%"),
    @("D49", "0.004203"),
    @("E49", "2.52%"),
    @("E50", "-0.14%"),
    @("D51", "0.0002001"),
    @("E51", "-0.14%")
)

foreach ($item in $changes) {
    $addr = $item[0]
    $val = $item[1]
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Write-Host "Applied $($changes.Count) changes"
